$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.010.31"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "4.035.33"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "535.75"
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").Value = "149.89"
$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("D7").Value = "4.031.51"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "0.698"
$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("E11").Value = "  -2.08%  "

$ws.Range("D12").Value = "53.99"
$ws.Range("E12").Value = "  +7.09%  "

$ws.Range("D13").Value = "0.0000328"
$ws.Range("E13").Value = "  -2.45%  "

$ws.Range("E14").Value = "  -2.07%  "

$ws.Range("D15").Value = "4.672.14"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "4.030.29"
$ws.Range("E16").Value = "  -0.51%  "

$ws.Range("D17").Value = "14.22"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "20.88"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("E19").Value = "  -3.24%  "

$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").Value = "71.960.30"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").Value = "433.44"
$ws.Range("E22").Value = "  -1.33%  "

$ws.Range("D23").Value = "98.32"
$ws.Range("E23").Value = "  -2.60%  "

$ws.Range("D24").Value = "3.59"
$ws.Range("E24").Value = "  -2.59%  "

$ws.Range("E25").Value = "  -1.94%  "

$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  -0.87%  "

$ws.Range("D27").Value = "4.43"
$ws.Range("E27").Value = "  +29.09%  "

$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").Value = "10.81"
$ws.Range("E29").Value = "  -2.20%  "

$ws.Range("D30").Value = "5.96"
$ws.Range("E30").Value = "  +2.11%  "

$ws.Range("D31").Value = "37.13"
$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").Value = "8.27"
$ws.Range("E32").Value = "  +22.05%  "

$ws.Range("E33").Value = "  +2.12%  "

$ws.Range("D34").Value = "50.15"
$ws.Range("E34").Value = "  +17.21%  "

$ws.Range("E35").Value = "  -0.92%  "

$ws.Range("D36").Value = "678.71"
$ws.Range("E36").Value = "  +0.31%  "

$ws.Range("D37").Value = "68.20"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("E38").Value = "  +4.37%  "

$ws.Range("D39").Value = "0.0₃0825"
$ws.Range("E39").Value = "  -5.28%  "

$ws.Range("E40").Value = "  +8.41%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "3.42"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.149"
$ws.Range("E42").Value = "  -6.26%  "

$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").Value = "11.10"
$ws.Range("E44").Value = "  +15.45%  "

$ws.Range("D45").Value = "0.0494"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  -2.66%  "

$ws.Range("E48").Value = "  -3.93%  "

$ws.Range("D49").Value = "3.37"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "3.10"
$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").Value = "2.871.28"
$ws.Range("E51").Value = "  +9.96%  "
